$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 corrections ---
$ws.Range("X19").Value = 1
$ws.Range("AI19").Value = 0

# --- Row 30 ID correction ---
$ws.Range("A30").Value = "'0966669"

# --- New rows 31-39 ---
$ws.Range("A31").Value = "'6333333"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = 1
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("U31").Value = 1
$ws.Range("V31").Value = 1
$ws.Range("W31").Value = 0
$ws.Range("X31").Value = 0
$ws.Range("Y31").Value = 0
$ws.Range("Z31").Value = 0
$ws.Range("AA31").Value = 0
$ws.Range("AB31").Value = 0
$ws.Range("AC31").Value = 1
$ws.Range("AD31").Value = 0
$ws.Range("AE31").Value = 0
$ws.Range("AF31").Value = 1
$ws.Range("AG31").Value = 0
$ws.Range("AH31").Value = 1
$ws.Range("AI31").Value = 1
$ws.Range("AJ31").Value = 0
$ws.Range("AK31").Value = 0
$ws.Range("AL31").Value = 0
$ws.Range("AM31").Value = 0
$ws.Range("AN31").Value = 1
$ws.Range("AO31").Value = 0
$ws.Range("AP31").Value = 1
$ws.Range("AQ31").Value = 1
$ws.Range("AR31").Value = 1
$ws.Range("AS31").Value = 1
$ws.Range("AT31").Value = 0

$ws.Range("A32").Value = "'9600220"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 1
$ws.Range("O32").Value = 1
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = 1
$ws.Range("S32").Value = 0
$ws.Range("T32").Value = 0
$ws.Range("U32").Value = 1
$ws.Range("V32").Value = 0
$ws.Range("W32").Value = 0
$ws.Range("X32").Value = 1
$ws.Range("Y32").Value = 1
$ws.Range("Z32").Value = 1
$ws.Range("AA32").Value = 1
$ws.Range("AB32").Value = 1
$ws.Range("AC32").Value = 1
$ws.Range("AD32").Value = 0
$ws.Range("AE32").Value = 1
$ws.Range("AF32").Value = 1
$ws.Range("AG32").Value = 0
$ws.Range("AH32").Value = 1
$ws.Range("AI32").Value = 1
$ws.Range("AJ32").Value = 0
$ws.Range("AK32").Value = 1
$ws.Range("AL32").Value = 0
$ws.Range("AM32").Value = 0
$ws.Range("AN32").Value = 0
$ws.Range("AO32").Value = 0
$ws.Range("AP32").Value = 0
$ws.Range("AQ32").Value = 1
$ws.Range("AR32").Value = 0
$ws.Range("AS32").Value = 0
$ws.Range("AT32").Value = 0

$ws.Range("A33").Value = "'0090000"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 1
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 1
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 1
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 1
$ws.Range("S33").Value = 0
$ws.Range("T33").Value = 1
$ws.Range("U33").Value = 1
$ws.Range("V33").Value = 1
$ws.Range("W33").Value = 1
$ws.Range("X33").Value = 1
$ws.Range("Y33").Value = 1
$ws.Range("Z33").Value = 0
$ws.Range("AA33").Value = 0
$ws.Range("AB33").Value = 0
$ws.Range("AC33").Value = 1
$ws.Range("AD33").Value = 1
$ws.Range("AE33").Value = 1
$ws.Range("AF33").Value = 1
$ws.Range("AG33").Value = 0
$ws.Range("AH33").Value = 1
$ws.Range("AI33").Value = 0
$ws.Range("AJ33").Value = 1
$ws.Range("AK33").Value = 1
$ws.Range("AL33").Value = 0
$ws.Range("AM33").Value = 1
$ws.Range("AN33").Value = 1
$ws.Range("AO33").Value = 1
$ws.Range("AP33").Value = 1
$ws.Range("AQ33").Value = 1
$ws.Range("AR33").Value = 1
$ws.Range("AS33").Value = 0
$ws.Range("AT33").Value = 0

$ws.Range("A34").Value = "'2189"
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = 1
$ws.Range("K34").Value = 1
$ws.Range("L34").Value = 1
$ws.Range("M34").Value = 1
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0

$ws.Range("A35").Value = "'5555"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = 1
$ws.Range("K35").Value = 1
$ws.Range("L35").Value = 1
$ws.Range("M35").Value = 1
$ws.Range("N35").Value = 1
$ws.Range("O35").Value = 1
$ws.Range("P35").Value = 1
$ws.Range("Q35").Value = 1

$ws.Range("A36").Value = "'2222"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 1
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = 1
$ws.Range("K36").Value = 1
$ws.Range("L36").Value = 1
$ws.Range("M36").Value = 1
$ws.Range("N36").Value = 1
$ws.Range("O36").Value = 1
$ws.Range("P36").Value = 1
$ws.Range("Q36").Value = 1

$ws.Range("A37").Value = "088#"
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 1
$ws.Range("O37").Value = 1
$ws.Range("P37").Value = 1
$ws.Range("Q37").Value = 1

$ws.Range("A38").Value = "'0991"
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 1
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 1

$ws.Range("A39").Value = "'0234"
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 1
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 1
$ws.Range("K39").Value = 1
$ws.Range("L39").Value = 1
$ws.Range("M39").Value = 1
$ws.Range("N39").Value = 1
$ws.Range("O39").Value = 1
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0

